$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.831.44"
$ws.Range("E2").Value = "  -3.87%  "
$ws.Range("D3").Value = "'3.328.86"
$ws.Range("E3").Value = "  -5.49%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'183.10"
$ws.Range("E5").Value = "  -8.68%  "
$ws.Range("D6").Value = "'531.84"
$ws.Range("E6").Value = "  -3.36%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'3.324.09"
$ws.Range("E8").Value = "  -5.27%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("E10").Value = "  -4.86%  "
$ws.Range("E11").Value = "  -7.59%  "
$ws.Range("D12").Value = "'0.134"
$ws.Range("E12").Value = "  -6.09%  "
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("E14").Value = "  -6.79%  "
$ws.Range("D15").Value = "'3.845.86"
$ws.Range("E15").Value = "  -6.49%  "
$ws.Range("D16").Value = "'3.319.69"
$ws.Range("E16").Value = "  -6.39%  "
$ws.Range("D17").Value = "'0.118"
$ws.Range("E17").Value = "  -5.12%  "
$ws.Range("D18").Value = "'17.87"
$ws.Range("E18").Value = "  -3.86%  "
$ws.Range("D19").Value = "'64.647.70"
$ws.Range("E19").Value = "  -4.08%  "
$ws.Range("D20").Value = "'11.21"
$ws.Range("E20").Value = "  -5.12%  "
$ws.Range("E21").Value = "  -6.33%  "
$ws.Range("D22").Value = "'376.61"
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("D23").Value = "'3.85"
$ws.Range("E23").Value = "  -4.17%  "
$ws.Range("D24").Value = "'11.34"
$ws.Range("E24").Value = "  -5.05%  "
$ws.Range("D25").Value = "'81.37"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").Value = "'3.92"
$ws.Range("E26").Value = "  +4.55%  "
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "'2.71"
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("D29").Value = "'11.67"
$ws.Range("E29").Value = "  -4.57%  "
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("D31").Value = "'29.21"
$ws.Range("E31").Value = "  -5.36%  "
$ws.Range("D32").Value = "'6.86"
$ws.Range("E32").Value = "  -4.54%  "
$ws.Range("D33").Value = "'646.79"
$ws.Range("E33").Value = "  -6.01%  "
$ws.Range("D34").Value = "'11.39"
$ws.Range("E34").Value = "  -3.40%  "
$ws.Range("D35").Value = "'59.92"
$ws.Range("E35").Value = "  -6.27%  "
$ws.Range("D36").Value = "'0.106"
$ws.Range("E36").Value = "  -4.81%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'0.398"
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("D39").Value = "'37.03"
$ws.Range("E39").Value = "  -4.71%  "
$ws.Range("D40").Value = "'0.0₃0734"
$ws.Range("E40").Value = "  +7.40%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("D43").Value = "'2.915.72"
$ws.Range("E43").Value = "  -5.33%  "
$ws.Range("D44").Value = "'2.54"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").Value = "'2.72"
$ws.Range("E45").Value = "  -9.53%  "
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("E47").Value = "  +11.24%  "
$ws.Range("D48").Value = "'2.66"
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("D49").Value = "'2.60"
$ws.Range("E49").Value = "  -7.52%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "'2.98"
$ws.Range("E51").Value = "  +0.81%  "
